$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Copy()
$ws.Range("B21").PasteSpecial(-4122)
Write-Host "done"
